$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "menu" column (old column D). Everything to its right (parent,
# level1, level2 headers / their data, and the trailing unused column width
# definition) shifts one column to the left.
$ws.Range("D1").EntireColumn.Delete() | Out-Null

# Update the first data row: path now points straight at the article, and
# its type changed from "Press Release" to "Article".
# (Set B2 before A2 so new shared-string entries are appended in the same
# order the source workbook uses: "Article" then the new path.)
$ws.Range("B2").Value = "Article"
$ws.Range("A2").Value = "/about-cancer/coping/feelings"

# Column A was manually widened to fit the new, longer path value.
$ws.Range("A1").EntireColumn.ColumnWidth = 51.75

# Active cell/selection moved to A2.
$ws.Range("A2").Select() | Out-Null
